$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new contact row right after the last used row in column A
# (xlUp = -4162), mirroring the "add record" flow the Electron UI drives.
$lastRow = $ws.Cells(1048576, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "rign"
$ws.Cells.Item($newRow, 2).Value = "rimg@raidio.com"
